# BOM.xlsx restructuring: add vendor/pricing detail rows, new parts, and brackets.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11: ruthex threaded inserts now sourced from 3dJake (not McMaster Carr)
# ---------------------------------------------------------------------------
$ws.Range("H11").Value = "3dJake"
$ws.Range("I11").Value = "RUT-GE-4-40x57-001"
$ws.Range("J11").Value = "ruthex"
$ws.Range("K11").Value = "GE-4-40x57-001"

# ---------------------------------------------------------------------------
# Row 20: relay - fill in pricing / vendor details
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = "electronics components"
$ws.Range("E2").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 1.43
$ws.Range("F20").Value = 20
$ws.Range("G20").Formula = '=$E20*$F20'
$ws.Range("G2").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("H20").Value = "Digikey"
$ws.Range("I20").Value = "1835-1117-ND"
$ws.Range("J20").Value = "Comus International"
$ws.Range("K20").Value = "3570-1331-053"

# ---------------------------------------------------------------------------
# Row 21: Neopixel - add vendor / part number detail
# ---------------------------------------------------------------------------
$ws.Range("H21").Value = "Adafruit"
$ws.Range("I21").Value = 1559
$ws.Range("J21").Value = "Adafruit"
$ws.Range("K21").Value = 1559

# ---------------------------------------------------------------------------
# Row 22: Arduino UNO - add link, pricing and vendor detail
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = "https://www.digikey.com/en/products/detail/dfrobot/DFR0216/6579366"
$ws.Hyperlinks.Add($ws.Range("D22"), "https://www.digikey.com/en/products/detail/dfrobot/DFR0216/6579366")
$ws.Range("E2").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = 16.9
$ws.Range("F22").Value = 1
$ws.Range("G22").Formula = '=$E22*$F22'
$ws.Range("G21").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "Digikey"
$ws.Range("I22").Value = "1738-1228-ND"
$ws.Range("J22").Value = "DFRobot"
$ws.Range("K22").Value = "DFR0216"
$ws.Range("D2").Copy()
$ws.Range("D22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 23: 5V-12V boost - add link and pricing
# ---------------------------------------------------------------------------
$ws.Range("D23").Value = "https://www.amazon.com/DROK-Boost-Converter-Regulator-Length/dp/B09M3LMSS3/"
$ws.Hyperlinks.Add($ws.Range("D23"), "https://www.amazon.com/DROK-Boost-Converter-Regulator-Length/dp/B09M3LMSS3/")
$ws.Range("E2").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = 7.99
$ws.Range("F23").Value = 1
$ws.Range("G23").Formula = '=$E23*$F23'
$ws.Range("G21").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 24: 12V to HV boost - add link and pricing
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = "https://www.amazon.com/Voltage-Converter-Vintage-Indicator-80V-380V/dp/B09D93QNYK"
$ws.Hyperlinks.Add($ws.Range("D24"), "https://www.amazon.com/Voltage-Converter-Vintage-Indicator-80V-380V/dp/B09D93QNYK")
$ws.Range("E2").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E24").Value = 11.39
$ws.Range("F24").Value = 1
$ws.Range("G24").Formula = '=$E24*$F24'
$ws.Range("G21").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D24").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 25: USB A to USB B cable - add link, pricing and vendor detail
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = "https://www.digikey.com/en/products/detail/assmann-wsw-components/AK672-2-2/947492"
$ws.Hyperlinks.Add($ws.Range("D25"), "https://www.digikey.com/en/products/detail/assmann-wsw-components/AK672-2-2/947492")
$ws.Range("E2").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = 1.64
$ws.Range("F25").Value = 1
$ws.Range("G25").Formula = '=$E25*$F25'
$ws.Range("G21").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("H25").Value = "Digikey"
$ws.Range("I25").Value = "AE1493-ND"
$ws.Range("J25").Value = "Assmann"
$ws.Range("K25").Value = "AK672/2-2"
$ws.Range("D2").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 26: 2 wire header cable - add link and pricing
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = "https://www.amazon.com/ZYAMY-Dupont-Connector-Multicolor-Breadboard/dp/B0B8Z23NWX/"
$ws.Hyperlinks.Add($ws.Range("D26"), "https://www.amazon.com/ZYAMY-Dupont-Connector-Multicolor-Breadboard/dp/B0B8Z23NWX/")
$ws.Range("E2").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = 5.99
$ws.Range("F26").Value = 1
$ws.Range("G26").Formula = '=$E26*$F26'
$ws.Range("G21").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D26").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Rows 27-29: quantities for remaining wire/header rows
# ---------------------------------------------------------------------------
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = 1

# ---------------------------------------------------------------------------
# Rows 30-33: custom PCB pricing from PCBWay
# ---------------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = 5.98
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = "PCBWay"

$ws.Range("E2").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = 1.82
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = "PCBWay"

$ws.Range("E2").Copy()
$ws.Range("E32").PasteSpecial(-4122)
$ws.Range("E32").Value = 1.82
$ws.Range("F32").Value = 1
$ws.Range("H32").Value = "PCBWay"

$ws.Range("E2").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = 1.4
$ws.Range("F33").Value = 1
$ws.Range("H33").Value = "PCBWay"

# ---------------------------------------------------------------------------
# Rows 34-37: 3d print brackets - add GitHub STL links and quantities
# ---------------------------------------------------------------------------
$ws.Range("D34").Value = "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/arduino-bracket.STL"
$ws.Hyperlinks.Add($ws.Range("D34"), "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/arduino-bracket.STL")
$ws.Range("F34").Value = 1
$ws.Range("D2").Copy()
$ws.Range("D34").PasteSpecial(-4122)

$ws.Range("D35").Value = "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/HV-DCDC-bracket.STL"
$ws.Hyperlinks.Add($ws.Range("D35"), "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/HV-DCDC-bracket.STL")
$ws.Range("F35").Value = 1
$ws.Range("D2").Copy()
$ws.Range("D35").PasteSpecial(-4122)

$ws.Range("D36").Value = "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/front-panel-center-bracket.STL"
$ws.Hyperlinks.Add($ws.Range("D36"), "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/front-panel-center-bracket.STL")
$ws.Range("F36").Value = 2
$ws.Range("D2").Copy()
$ws.Range("D36").PasteSpecial(-4122)

$ws.Range("D37").Value = "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/front-panel-edge-bracket.STL"
$ws.Hyperlinks.Add($ws.Range("D37"), "https://github.com/lafefspietz/MEMSduino/blob/main/3dprint_files/front-panel-edge-bracket.STL")
$ws.Range("F37").Value = 2
$ws.Range("D2").Copy()
$ws.Range("D37").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 38: DB25 to Fischer Cable - complete row
# ---------------------------------------------------------------------------
$ws.Range("B38").Value = "system"
$ws.Range("B29").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("C38").Value = "wires"
$ws.Range("D38").Value = "https://aivon.fi"
$ws.Hyperlinks.Add($ws.Range("D38"), "https://aivon.fi")
$ws.Range("F38").Value = 1
$ws.Range("D2").Copy()
$ws.Range("D38").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Rows 39-42: new machined brass H Bracket hardware rows
# ---------------------------------------------------------------------------
$ws.Range("A39").Value = "H Bracket"
$ws.Range("B39").Value = "cold assembly"
$ws.Range("B2").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$ws.Range("C39").Value = "machined brass"

$ws.Range("B40").Value = "cold assembly"
$ws.Range("B2").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("C40").Value = "machined brass"

$ws.Range("B41").Value = "cold assembly"
$ws.Range("B2").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("C41").Value = "machined brass"

$ws.Range("B42").Value = "cold assembly"
$ws.Range("B2").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("C42").Value = "machined brass"

# ---------------------------------------------------------------------------
# View tidy-up: active cell moved back to D4 (no frozen top-left scroll)
# ---------------------------------------------------------------------------
$ws.Range("D4").Select()
